# Updated model coefficients for "Common hardwoods" (new Fire:Size / Fire:Northness
# interaction model), plus a couple of small "Fate" index corrections elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Fate" column header (H) label removed on the two tables that still had
#    the stray header text ("Fate") above the per-row model-id column; the
#    cell itself stays (with the same numeric format as the rest of the
#    coefficient block) but now carries no text.
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = ""
$ws.Range("H2").NumberFormat = "0.00"
$ws.Range("H6").Value = ""
$ws.Range("H6").NumberFormat = "0.00"

# H9 / H13 previously had no cell at all under the other two tables; give them
# the same (empty, numeric-formatted) treatment so the column is consistent.
$ws.Range("H9").NumberFormat = "0.00"
$ws.Range("H13").NumberFormat = "0.00"

# The whole "Fate id" column (H) picks up the 0.00 number format used by the
# coefficient cells next to it.
$ws.Range("H3").NumberFormat = "0.00"
$ws.Range("H4").NumberFormat = "0.00"
$ws.Range("H7").NumberFormat = "0.00"
$ws.Range("H10").NumberFormat = "0.00"
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H14").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 2. "Common hardwoods" interactions label + new interaction coefficient
#    columns (Q:AB) with headers in row 2.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "Fire:Size; Fire:Northness"

$ws.Range("Q2").Value = "fsCat1:northness"
$ws.Range("R2").Value = "fsCat2:northness"
$ws.Range("S2").Value = "fsCat3:northness"
$ws.Range("T2").Value = "TSizeCatTR1:fsCat1"
$ws.Range("U2").Value = "TSizeCatTR2:fsCat1"
$ws.Range("V2").Value = "TSizeCatTR3:fsCat1"
$ws.Range("W2").Value = "TSizeCatTR1:fsCat2"
$ws.Range("X2").Value = "TSizeCatTR2:fsCat2"
$ws.Range("Y2").Value = "TSizeCatTR3:fsCat2"
$ws.Range("Z2").Value = "TSizeCatTR1:fsCat3"
$ws.Range("AA2").Value = "TSizeCatTR2:fsCat3"
$ws.Range("AB2").Value = "TSizeCatTR3:fsCat3"
$ws.Range("Q2:AB2").NumberFormat = "0.00"

# Updated "Common hardwoods" coefficients, Fate 1 (row 3) and Fate 2 (row 4),
# including the 12 new interaction-term columns (Q:AB).
$ws.Range("I3").Value = -48.745599769925597
$ws.Range("J3").Value = -11.9116843011664
$ws.Range("K3").Value = -28.111517637552598
$ws.Range("L3").Value = 2.4617176961584302
$ws.Range("M3").Value = 49.121575488957497
$ws.Range("N3").Value = 48.662570208942
$ws.Range("O3").Value = 48.968179831545399
$ws.Range("P3").Value = 113.149905651034
$ws.Range("Q3").Value = -115.37580436776901
$ws.Range("R3").Value = -112.732880524892
$ws.Range("S3").Value = -113.430889607518
$ws.Range("T3").Value = 12.3968233494017
$ws.Range("U3").Value = 26.286124676134101
$ws.Range("V3").Value = -4.6184443294156399
$ws.Range("W3").Value = 11.932848627102301
$ws.Range("X3").Value = 28.257116276785499
$ws.Range("Y3").Value = -3.4364280206032101
$ws.Range("Z3").Value = 12.651242700557299
$ws.Range("AA3").Value = 29.4046904535987
$ws.Range("AB3").Value = -1.3739425637656699

$ws.Range("I4").Value = 3.5886078648774502
$ws.Range("J4").Value = 0.40963981415978101
$ws.Range("K4").Value = -1.43522993973628
$ws.Range("L4").Value = -0.55337346475902405
$ws.Range("M4").Value = -5.5793283272545402
$ws.Range("N4").Value = -78.283956735203205
$ws.Range("O4").Value = -93.330840633817701
$ws.Range("P4").Value = -1.5972560547421299
$ws.Range("Q4").Value = 0.88597935747399204
$ws.Range("R4").Value = 3.1772367954099199
$ws.Range("S4").Value = 1.14593496669057
$ws.Range("T4").Value = 1.41616236689124
$ws.Range("U4").Value = 5.4505343206587904
$ws.Range("V4").Value = 5.5066099172658696
$ws.Range("W4").Value = 71.4584686551148
$ws.Range("X4").Value = 77.430876355669596
$ws.Range("Y4").Value = 77.560921334759399
$ws.Range("Z4").Value = 30.793697972078199
$ws.Range("AA4").Value = 89.070122043649306
$ws.Range("AB4").Value = 91.434914507301698

# ---------------------------------------------------------------------------
# 3. Two "Fate" id corrections (PSEMEN row was mislabeled 2, ARCMAN row was
#    mislabeled 2; both are the first/only model instance so should read 1).
# ---------------------------------------------------------------------------
$ws.Range("H7").Value = 1
$ws.Range("H14").Value = 1

# ---------------------------------------------------------------------------
# 4. Column F is now wide enough to show the new "Fire:Size; Fire:Northness"
#    label in full.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 20.6666667

# ---------------------------------------------------------------------------
# 5. Selection / view state left pointing at the block that was edited.
# ---------------------------------------------------------------------------
$ws.Range("H9:N11").Select()
